# Plantilla_estandar_para_cargar_servicios.xlsx
# "se modifican plantillas de excel"
#
# Summary of the edit:
#  - Hoja1!F2 / G2 headers get "/ ID de Servicios" appended.
#  - Desplegables!B4 gains a new entry "Id de Servicios" (column B now has
#    4 entries instead of 3).
#  - Data validation on Hoja1!F3 becomes its own rule allowing the 4-item
#    list (Desplegables!$B$2:$B$4); F4:F102 keeps the original 3-item list.
#  - Columns F and G on Hoja1 get wider to fit the longer header text.
#  - The active selection moves to G2.

$wb = $excel.ActiveWorkbook

$hoja1 = $wb.Worksheets.Item("Hoja1")
$desplegables = $wb.Worksheets.Item("Desplegables")

# --- Hoja1 header text + Desplegables new option (order matters for the
#     shared-strings table layout, so mirror the original authoring order:
#     F2 header, then the new dropdown option, then the G2 header) ----------
$hoja1.Range("F2").Value = "Centro de costos/ Orden de inversión/ ID de Servicios"

# Desplegables column B previously listed 3 options (rows 1-3). Add a 4th.
$desplegables.Range("B4").Value = "Id de Servicios"

$hoja1.Range("G2").Value = "Número centro de costos/ Orden de inversión/ ID de Servicios"

# --- Column widths to fit the new, longer header text -----------------------
$hoja1.Columns.Item(6).ColumnWidth = 61
$hoja1.Columns.Item(7).ColumnWidth = 68.1

# --- Data validation: split F3:F102 into F3 (4-item list) and F4:F102 (3-item list)
$hoja1.Range("F3:F102").Validation.Delete()

$vF4 = $hoja1.Range("F4:F102").Validation
$vF4.Add(3, 1, 1, "=Desplegables!`$B`$2:`$B`$3")
$vF4.InputMessage = "Seleccione"

$vF3 = $hoja1.Range("F3").Validation
$vF3.Add(3, 1, 1, "=Desplegables!`$B`$2:`$B`$4")
$vF3.InputMessage = "Seleccione"

# --- Move the active selection to G2 ----------------------------------------
$hoja1.Activate() | Out-Null
$hoja1.Range("G2").Select() | Out-Null
